$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.9969033162691203
$ws.Range("C3").Value = 0.9972067099356036
$ws.Range("D3").Value = 0.9971279349314993

# Row 4 - GradientBoostingRegressor
$ws.Range("B4").Value = 0.9965699670864142
$ws.Range("C4").Value = 0.9965700457374805
$ws.Range("D4").Value = 0.9965699670864142

# Row 5 - AdaBoostRegressor
$ws.Range("B5").Value = 0.985639071626108
$ws.Range("C5").Value = 0.9854972758612843
$ws.Range("D5").Value = 0.9861078689948336
